# Generate Report for Handoff
# Replace the old GUID-based filenames / timestamps with the new handoff run's
# values, and clear the "Latest Target File" / "Latest Handback File" columns
# (handback hasn't happened yet for this run).

$wb = $excel.ActiveWorkbook

$oldGuid = "edfe7190-de0f-42f5-a0a9-6e75d6b9ddb8"
$newGuid = "9c9b7c0d-7d7b-4729-9fa3-d711e3b35526"

$oldXlfBase = "3cc98770d028d9db8104a0e1657d3075796f3342"
$newXlfBase = "3229d8d6d4b73bfabea289e139145ea78a2120a9"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newXlfBase.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-18 02:52:38"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsZhCn.Hyperlinks.Item(2).Delete()

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newXlfBase.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-18 02:52:43"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDeDe.Hyperlinks.Item(2).Delete()
